$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.592.19"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").Value = "2.528.89"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'316.28"
$ws.Range("E5").Value = "  +3.27%  "

$ws.Range("D6").Value = "'95.13"
$ws.Range("E6").Value = "  -5.63%  "

$ws.Range("E7").Value = "  -1.17%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  -4.10%  "

$ws.Range("D10").Value = "'36.18"
$ws.Range("E10").Value = "  -3.80%  "

$ws.Range("D11").Value = "'0.0807"
$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'7.54"
$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.113"
$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").Value = "2.919.78"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.570.25"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'15.45"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").Value = "'0.847"
$ws.Range("E17").Value = "  -2.51%  "

$ws.Range("D18").Value = "42.628.63"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("D19").Value = "'12.97"
$ws.Range("E19").Value = "  -1.46%  "

$ws.Range("D20").Value = "'6.56"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -2.93%  "

$ws.Range("D22").Value = "'70.04"
$ws.Range("E22").Value = "  -2.39%  "

$ws.Range("D23").Value = "'250.67"
$ws.Range("E23").Value = "  -1.58%  "

$ws.Range("D24").Value = "'2.96"
$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("E25").Value = "  -2.58%  "

$ws.Range("D26").Value = "'26.58"
$ws.Range("E26").Value = "  -3.33%  "

$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").Value = "'2.39"
$ws.Range("E28").Value = "  +2.01%  "

$ws.Range("D29").Value = "'39.15"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D30").Value = "'10.15"
$ws.Range("E30").Value = "  -3.54%  "

$ws.Range("D31").Value = "'6.05"
$ws.Range("E31").Value = "  -1.97%  "

$ws.Range("D32").Value = "'155.08"
$ws.Range("E32").Value = "  -2.21%  "

$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("D34").Value = "'19.08"
$ws.Range("E34").Value = "  +3.00%  "

$ws.Range("D35").Value = "'3.28"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("D36").Value = "'0.0784"
$ws.Range("E36").Value = "  -2.23%  "

$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("E38").Value = "  -4.36%  "

$ws.Range("D39").Value = "'0.118"
$ws.Range("E39").Value = "  -1.25%  "

$ws.Range("D40").Value = "'23.61"
$ws.Range("E40").Value = "  -2.25%  "

$ws.Range("E41").Value = "  +11.41%  "

$ws.Range("D42").Value = "'3.80"
$ws.Range("E42").Value = "  -3.31%  "

$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("D44").Value = "'0.0300"
$ws.Range("E44").Value = "  -1.64%  "

$ws.Range("E45").Value = "  -5.73%  "

$ws.Range("D46").Value = "2.017.39"
$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("D47").Value = "'85.44"
$ws.Range("E47").Value = "  -0.94%  "

$ws.Range("D48").Value = "'8.79"
$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("D49").Value = "2.775.05"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("D50").Value = "'73.80"
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("D51").Value = "'102.47"
$ws.Range("E51").Value = "  -1.39%  "
